# edit.ps1 -- "Updated symbol list on Sun Jan 15 08:51:05 UTC 2023 with GitHub Actions"
#
# The crypto-tracker sheet refreshes its Price (D) / Volume(1h) (E) columns on
# every scheduled run, and (for this run) two rows swap which coin currently
# occupies rank #6 / #7 (FTXToken <-> GateToken, including their Link cells).
#
# Every touched D/E cell in the source workbook is stored as literal TEXT
# (e.g. "7.720", "-5.58%"), not a Number/Percentage, so formatting/precision
# has to be preserved exactly. A plain `.Value = "7.720"` on a General-format
# cell would make Excel "helpfully" reinterpret it as a number and silently
# drop the trailing zero (or turn "-5.58%" into a 0.0558-style fraction), so
# each target cell is first marked as Text ("@") -- exactly what a human author
# would do via Format Cells before typing these values in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "295.97"
$ws.Range("E2").Value = "-5.58%"

# Row 3: OKB
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "31.46"
$ws.Range("E3").Value = "-3.00%"

# Row 4: HuobiToken
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.111"
$ws.Range("E4").Value = "-4.17%"

# Row 5: Cronos
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07463"
$ws.Range("E5").Value = "-2.92%"

# Row 6: KuCoinToken
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "7.720"
$ws.Range("E6").Value = "-2.02%"

# Row 7: FTXToken -> GateToken
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "3.785"
$ws.Range("E7").Value = "1.91%"

# Row 8: GateToken -> FTXToken
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "1.687"
$ws.Range("E8").Value = "3.61%"

# Row 9: MXToken
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9297"
$ws.Range("E9").Value = "1.12%"

# Row 10: WazirX
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1684"
$ws.Range("E10").Value = "-2.46%"

# Row 11: LiechtensteinCryptoassetsExchange
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07184"
$ws.Range("E11").Value = "-5.42%"

# Row 12: MandalaExchangeToken
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07884"
$ws.Range("E12").Value = "-4.30%"

# Row 13: BitrueCoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03011"
$ws.Range("E13").Value = "-0.74%"

# Row 14: BitMartToken
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09919"
$ws.Range("E14").Value = "0.28%"

# Row 15: BitForexToken
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").Value = "-1.04%"

# Row 16: TigerCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006278"
$ws.Range("E16").Value = "2.98%"

# Row 17: LEO
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.457"
$ws.Range("E17").Value = "-0.50%"

# Row 18: BTSEToken
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "2.222"
$ws.Range("E18").Value = "-1.04%"

# Row 19: BitpandaEcosystemToken
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3280"
$ws.Range("E19").Value = "-1.32%"

# Row 20: ProBitToken
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1336"
$ws.Range("E20").Value = "0.09%"

# Row 21: MCDex
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "4.576"
$ws.Range("E21").Value = "7.62%"

# Row 22: CoinExToken
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04649"
$ws.Range("E22").Value = "1.83%"

# Row 23: ZBToken
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1553"
$ws.Range("E23").Value = "-4.38%"

# Row 24: BitKan
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "-0.22%"

# Row 25: HotbitToken
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.99%"

# Row 26: NitroEx
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "0.36%"

# Row 27: UpBots
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001878"
$ws.Range("E27").Value = "8.11%"

# Row 39: One
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01657"
$ws.Range("E39").Value = "-6.34%"

# Row 40: IDEX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04435"
$ws.Range("E40").Value = "-4.64%"

# Row 41: KickToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007096"
$ws.Range("E41").Value = "-1.38%"

# Row 42: BKEXToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").Value = "-3.30%"

# Row 43: CEJI
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002073"
$ws.Range("E43").Value = "-8.07%"

# Row 44: LocalTraders
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-16.15%"

# Row 45: CoinLion
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006014"
$ws.Range("E45").Value = "-3.04%"

# Row 46: BOLO
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = "1.918"
$ws.Range("E46").Value = "1.34%"

# Row 47: CoinbaseStockToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01101"
$ws.Range("E47").Value = "-15.09%"
